$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.817.88"
$ws.Range("E2").Value = "  +4.87%  "
$ws.Range("D3").Value = "2.353.42"
$ws.Range("E3").Value = "  +4.69%  "
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.03"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.08"
$ws.Range("E6").Value = "  +4.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +4.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.08"
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.47"
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "2.710.97"
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("D15").Value = "2.350.07"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.33"
$ws.Range("E16").Value = "  +5.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.835"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "46.774.81"
$ws.Range("E18").Value = "  +5.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("E19").Value = "  +15.35%  "
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.90"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.86"
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.94"
$ws.Range("E27").Value = "  +13.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.93"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.31"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.77"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.13"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0821"
$ws.Range("E33").Value = "  +4.78%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.15"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.07"
$ws.Range("E39").Value = "  +7.95%  "
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.03"
$ws.Range("E42").Value = "  -7.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  +10.43%  "
$ws.Range("D45").Value = "1.812.49"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.30"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.85"
$ws.Range("E48").Value = "  +7.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.94"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.11"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.65"
$ws.Range("E51").Value = "  +3.34%  "
